$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Add the Organizations list on Feuil1!J3:J6 ---------------------------
$ws.Range("J3").Value = "ATI"
$ws.Range("J4").Value = "LABEX"
$ws.Range("J5").Value = "TeamEd"
$ws.Range("J6").Value = "Other"

# --- Defined name used by the dropdown list --------------------------------
$wb.Names.Add("Organizations", "=Feuil1!`$J`$3:`$J`$6")

# --- Data validation (dropdown list) on G2 ---------------------------------
$ws.Range("G2").Validation.Add(3, 1, 1, "=Organizations")
$ws.Range("G2").Validation.IgnoreBlank = $false

# --- Column widths (header row formatting) ---------------------------------
$ws.Columns.Item(1).ColumnWidth = 13.022135416666666
$ws.Columns.Item(2).ColumnWidth = 13.307291666666666
$ws.Columns.Item(3).ColumnWidth = 16.592447916666668
$ws.Columns.Item(4).ColumnWidth = 14.451822916666666
$ws.Columns.Item(5).ColumnWidth = 15.022135416666666
$ws.Columns.Item(6).ColumnWidth = 12.736979166666666
$ws.Columns.Item(7).ColumnWidth = 14.166666666666666

# --- Shrink the XML-mapped table back down to the actual data range -------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:G2")) | Out-Null

# --- Page setup (adds <pageSetup>) -----------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection moves to F2 --------------------------------------------------
$ws.Range("F2").Select() | Out-Null
